# Weekly update: a new price-observation date is added for "Cilantro" at
# "Vega Monumental Concepción", and the oldest existing observation date
# drops off the rolling window (duplicated at the bottom, matching the
# historical pattern already present in the sheet).
#
# Concretely: two new rows are inserted right after the header block (at
# row 86, before the first data row of this rolling window), pushing the
# existing rows 86-153 down to 88-155. The two freshly inserted rows
# (86-87) are populated with the same "Primera"/"Segunda" quality-pair
# template as the rest of the sheet, but stamped with the new date
# (2022-01-06). The two rows that fall off the end of the original range
# (what were rows 152-153) reappear unchanged as the new rows 154-155.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstDataRow = 86

# Snapshot the two template rows (A:R) *before* inserting anything, since
# Range objects are resolved by address, not tracked through the insert.
$templatePrimeraRow = $ws.Range("A" + $firstDataRow + ":R" + $firstDataRow).Value()
$templateSegundaRow = $ws.Range("A" + ($firstDataRow + 1) + ":R" + ($firstDataRow + 1)).Value()

# Insert two blank rows at the top of the data window; this shifts the
# existing rows 86-153 down to 88-155 (and copies row formatting/number
# formats down along with them), and grows the sheet dimension to R155.
$ws.Range("A" + $firstDataRow + ":R" + ($firstDataRow + 1)).EntireRow.Insert()

# Re-populate the newly inserted (now-blank) rows with the template data.
$ws.Range("A" + $firstDataRow + ":R" + $firstDataRow).Value = $templatePrimeraRow
$ws.Range("A" + ($firstDataRow + 1) + ":R" + ($firstDataRow + 1)).Value = $templateSegundaRow

# Stamp the new rows with this week's date (2022-01-06), replacing the
# date copied in from the template.
$newDateSerial = 44567
$ws.Range("D" + $firstDataRow).Value = $newDateSerial
$ws.Range("D" + ($firstDataRow + 1)).Value = $newDateSerial
